$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.048.68"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").Value = "3.090.32"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.42%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.088.78"
$ws.Range("E8").Value = "  +0.54%  "

$ws.Range("E9").Value = "  +3.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.05%  "

$ws.Range("E11").Value = "  -0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.401"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.44%  "

$ws.Range("D13").Value = "3.624.42"
$ws.Range("E13").Value = "  +0.73%  "

$ws.Range("E14").Value = "  +1.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "57.173.20"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").Value = "3.093.01"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "349.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "

$ws.Range("E26").Value = "  -1.40%  "

$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("D29").Value = "0.0₃0883"
$ws.Range("E29").Value = "  -2.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("E31").Value = "  +1.83%  "

$ws.Range("E32").Value = "  +0.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  +9.24%  "

$ws.Range("E36").Value = "  -2.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.10%  "

$ws.Range("E38").Value = "  -1.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  -0.65%  "

$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("E43").Value = "  +1.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.694"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.26%  "

$ws.Range("D45").Value = "2.387.39"
$ws.Range("E45").Value = "  +5.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").Value = "3.130.29"
$ws.Range("E48").Value = "  +0.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0264"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.40%  "

$ws.Range("E50").Value = "  -2.27%  "

$ws.Range("E51").Value = "  -1.17%  "
